# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# described in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    # Force the cell to be treated as text so that numeric-looking
    # strings such as "4.03" are not coerced into numbers, then
    # restore the default "Normal" style so no visible formatting
    # change is introduced (matches the original plain inline-string cells).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '26.268.16'
$ws.Range('E2').Value = '  +0.52%  '
Set-TextCell 'D3' '1.592.70'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  -0.26%  '
Set-TextCell 'D5' '212.75'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('E9').Value = '  +0.03%  '
Set-TextCell 'D11' '0.0850'
$ws.Range('E11').Value = '  +0.62%  '
Set-TextCell 'D12' '1.817.71'
$ws.Range('E12').Value = '  +1.00%  '
Set-TextCell 'D13' '1.595.01'
$ws.Range('E13').Value = '  -0.49%  '
Set-TextCell 'D14' '4.03'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('E15').Value = '  +1.59%  '
Set-TextCell 'D16' '64.37'
$ws.Range('E16').Value = '  -0.09%  '
Set-TextCell 'D17' '26.274.47'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('E18').Value = '  -0.45%  '
Set-TextCell 'D19' '7.43'
$ws.Range('E19').Value = '  +2.53%  '
Set-TextCell 'D20' '213.30'
$ws.Range('E20').Value = '  +2.69%  '
$ws.Range('E21').Value = '  -0.22%  '
Set-TextCell 'D22' '4.28'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('E23').Value = '  +1.68%  '
Set-TextCell 'D24' '2.14'
$ws.Range('E24').Value = '  -3.14%  '
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  -0.21%  '
Set-TextCell 'D27' '7.07'
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('E28').Value = '  -0.48%  '
Set-TextCell 'D29' '15.19'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('E31').Value = '  +1.17%  '
Set-TextCell 'D32' '3.20'
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('E33').Value = '  -0.39%  '
Set-TextCell 'D34' '1.338.06'
$ws.Range('E34').Value = '  +4.76%  '
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('E36').Value = '  -0.45%  '
Set-TextCell 'D37' '0.591'
$ws.Range('E37').Value = '  -2.94%  '
Set-TextCell 'D38' '0.0166'
$ws.Range('E38').Value = '  +0.32%  '
Set-TextCell 'D39' '0.821'
$ws.Range('E39').Value = '  +0.73%  '
Set-TextCell 'D40' '1.03'
$ws.Range('E40').Value = '  -9.14%  '
Set-TextCell 'D41' '5.73'
$ws.Range('E41').Value = '  +3.11%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('E43').Value = '  +0.17%  '
Set-TextCell 'D44' '0.765'
$ws.Range('E44').Value = '  +0.17%  '
Set-TextCell 'D45' '62.02'
$ws.Range('E45').Value = '  -0.67%  '
Set-TextCell 'D46' '1.729.73'
$ws.Range('E46').Value = '  +0.89%  '
Set-TextCell 'D47' '85.92'
$ws.Range('E47').Value = '  -3.29%  '
$ws.Range('E48').Value = '  -4.06%  '
Set-TextCell 'D50' '0.0975'
$ws.Range('E50').Value = '  -2.92%  '
Set-TextCell 'D51' '0.999'
$ws.Range('E51').Value = '  -0.28%  '
